# TC01_CDS_Filter_Acesses-Controlled.xlsx edit
# Commit: "Created CDS TC01,2 for Acesses"
#
# The three data rows (2-4) on the "startup" sheet had their long Neo4j
# Cypher query text (columns B "WebExcel query" and C "dbExcel query")
# cleared out, while keeping the existing cell formatting/style. The
# referenced "StatQuery" filenames in columns D/E are unchanged values.
# The active cell selection on the sheet moved from C4 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Clear the long query text from columns B and C for rows 2-4, keeping
# the existing cell style (wrap-text format) intact.
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""

# The custom (very tall) row heights that were sized for the long query
# text are no longer needed now that the cells are empty - let Excel
# recompute the natural row height.
$ws.Rows.Item(2).AutoFit() | Out-Null
$ws.Rows.Item(3).AutoFit() | Out-Null
$ws.Rows.Item(4).AutoFit() | Out-Null

# Update the sheet's active selection to A2 (was C4).
$ws.Range("A2").Select() | Out-Null
